$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Japanese_Yen nominal changed from 1000 to 100 (mr_parser addition of tenge)
# and the corresponding converted value recalculated accordingly.

$ws.Range("D8").Value = 100
$ws.Range("E8").Value = 52.5113

$ws.Range("D16").Value = 100
$ws.Range("E16").Value = 53.2682

$ws.Range("D24").Value = 100
$ws.Range("E24").Value = 53.7028

$ws.Range("D32").Value = 100
$ws.Range("E32").Value = 53.0609
